$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number and date range) ---
$ws.Range("A8").Value = "Volume 32   Number  31"
$ws.Range("C9").Value = "Report Covering the Week  7/28/2025  Through  8/3/2025"

# --- Weekly crime-stat table cell updates ---
$ws.Range("F15").NumberFormat = "@"
$ws.Range("F15").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("F15").PasteSpecial(-4122)
$ws.Range("G15").NumberFormat = "@"
$ws.Range("G15").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("G15").PasteSpecial(-4122)
$ws.Range("H15").NumberFormat = "@"
$ws.Range("H15").Value = "***.*"
$ws.Range("C14").Copy()
$ws.Range("H15").PasteSpecial(-4122)
$ws.Range("N15").Value = -50
$ws.Range("C16").Value = 1
$ws.Range("D16").Value = 2
$ws.Range("E16").Value = -50
$ws.Range("F16").Value = 6
$ws.Range("G16").Value = 4
$ws.Range("H16").Value = 50
$ws.Range("I16").Value = 58
$ws.Range("J16").Value = 73
$ws.Range("K16").Value = -20.547945205479
$ws.Range("L16").Value = -9.375
$ws.Range("M16").Value = -57.352941176470
$ws.Range("N16").Value = -85.888077858880
$ws.Range("C17").Value = 4
$ws.Range("D17").Value = 2
$ws.Range("E17").Value = 100
$ws.Range("F17").Value = 13
$ws.Range("G17").Value = 16
$ws.Range("H17").Value = -18.75
$ws.Range("I17").Value = 89
$ws.Range("J17").Value = 96
$ws.Range("K17").Value = -7.291666666666
$ws.Range("L17").Value = -10.101010101010
$ws.Range("M17").Value = 27.142857142857
$ws.Range("N17").Value = -53.403141361256
$ws.Range("C18").Value = 1
$ws.Range("I14").Copy()
$ws.Range("C18").PasteSpecial(-4122)
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 2
$ws.Range("G18").Value = 6
$ws.Range("H18").Value = -66.666666666666
$ws.Range("I18").Value = 38
$ws.Range("J18").Value = 39
$ws.Range("K18").Value = -2.564102564102
$ws.Range("L18").Value = -29.629629629629
$ws.Range("M18").Value = -74.496644295302
$ws.Range("N18").Value = -94.594594594594
$ws.Range("C19").Value = 11
$ws.Range("D19").Value = 7
$ws.Range("E19").Value = 57.142857142857
$ws.Range("F19").Value = 37
$ws.Range("G19").Value = 31
$ws.Range("H19").Value = 19.354838709677
$ws.Range("I19").Value = 299
$ws.Range("J19").Value = 337
$ws.Range("K19").Value = -11.275964391691
$ws.Range("L19").Value = -22.337662337662
$ws.Range("M19").Value = 8.333333333333
$ws.Range("N19").Value = -23.333333333333
$ws.Range("D20").Value = 7
$ws.Range("E20").Value = -42.857142857142
$ws.Range("F20").Value = 14
$ws.Range("G20").Value = 22
$ws.Range("H20").Value = -36.363636363636
$ws.Range("I20").Value = 85
$ws.Range("J20").Value = 102
$ws.Range("K20").Value = -16.666666666666
$ws.Range("L20").Value = 13.333333333333
$ws.Range("M20").Value = -13.265306122449
$ws.Range("N20").Value = -94.842233009708
$ws.Range("C21").Value = 21
$ws.Range("D21").Value = 19
$ws.Range("E21").Value = 10.526315789473
$ws.Range("F21").Value = 72
$ws.Range("G21").Value = 79
$ws.Range("H21").Value = -8.860759493670
$ws.Range("I21").Value = 579
$ws.Range("J21").Value = 653
$ws.Range("K21").Value = -11.332312404287
$ws.Range("L21").Value = -15.597667638484
$ws.Range("M21").Value = -21.544715447154
$ws.Range("N21").Value = -82.829181494661
$ws.Range("C23").NumberFormat = "@"
$ws.Range("C23").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("C23").PasteSpecial(-4122)
$ws.Range("D23").Value = 1
$ws.Range("I14").Copy()
$ws.Range("D23").PasteSpecial(-4122)
$ws.Range("E23").Value = -100
$ws.Range("L14").Copy()
$ws.Range("E23").PasteSpecial(-4122)
$ws.Range("J23").Value = 22
$ws.Range("K23").Value = -45.454545454545
$ws.Range("C24").Value = 18
$ws.Range("E24").Value = -21.739130434782
$ws.Range("F24").Value = 81
$ws.Range("G24").Value = 91
$ws.Range("H24").Value = -10.989010989011
$ws.Range("I24").Value = 740
$ws.Range("J24").Value = 930
$ws.Range("K24").Value = -20.430107526881
$ws.Range("L24").Value = 2.635228848821
$ws.Range("M24").Value = 34.545454545454
$ws.Range("C25").Value = 11
$ws.Range("D25").Value = 23
$ws.Range("E25").Value = -52.173913043478
$ws.Range("F25").Value = 54
$ws.Range("G25").Value = 76
$ws.Range("H25").Value = -28.947368421052
$ws.Range("I25").Value = 522
$ws.Range("J25").Value = 771
$ws.Range("K25").Value = -32.295719844358
$ws.Range("L25").Value = 3.571428571428
$ws.Range("F26").Value = 26
$ws.Range("G26").Value = 19
$ws.Range("H26").Value = 36.842105263157
$ws.Range("I26").Value = 178
$ws.Range("J26").Value = 190
$ws.Range("K26").Value = -6.315789473684
$ws.Range("L26").Value = 11.25
$ws.Range("M26").Value = -14.832535885167
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("D27").PasteSpecial(-4122)
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "***.*"
$ws.Range("C14").Copy()
$ws.Range("E27").PasteSpecial(-4122)
$ws.Range("F27").NumberFormat = "@"
$ws.Range("F27").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("F27").PasteSpecial(-4122)
$ws.Range("G27").Value = 1
$ws.Range("H27").Value = -100
$ws.Range("C28").Value = 2
$ws.Range("I14").Copy()
$ws.Range("C28").PasteSpecial(-4122)
$ws.Range("D28").Value = 2
$ws.Range("E28").Value = 0
$ws.Range("F28").Value = 5
$ws.Range("H28").Value = -16.666666666666
$ws.Range("I28").Value = 18
$ws.Range("J28").Value = 21
$ws.Range("K28").Value = -14.285714285714
$ws.Range("L28").Value = 0
$ws.Range("I33").Value = 2
